$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Praktikumsabschnitt" in column H, matching style of existing headers (A1:G1)
$ws.Range("H1").Value = "Praktikumsabschnitt"
$ws.Range("H1").HorizontalAlignment = -4108   # xlCenter, same as other header cells

# Add the data value for the new column in row 2
$ws.Range("H2").Value = 1

# Set column H width to match the new content (closest achievable to 23.6328125)
$ws.Columns.Item(8).ColumnWidth = 22.8

# Update the active selection to reflect the newly used cell
$ws.Range("H6").Select()
